$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.920.61'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '3.776.69'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.46%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("E9").Value = '  -1.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.446'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.53'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.64%  '
$ws.Range("E12").Value = '  -2.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.29'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.68%  '
$ws.Range("D14").Value = '4.411.12'
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").Value = '3.775.05'
$ws.Range("E15").Value = '  -0.13%  '
$ws.Range("D16").Value = '67.928.23'
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.09%  '
$ws.Range("E18").Value = '  +2.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.99'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '458.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.694'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '82.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.82%  '
$ws.Range("E24").Value = '  -6.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.21%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.89'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.19%  '
$ws.Range("D29").Value = '3.929.02'
$ws.Range("E29").Value = '  -0.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("E32").Value = '  -7.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.98'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.01%  '
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.92'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0989'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("E37").Value = '  +2.55%  '
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.981'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.83%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.20'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.39'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.84%  '
$ws.Range("E44").Value = '  -1.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '152.92'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.295'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.77%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.06%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.84'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '387.24'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.50'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.67%  '
